$wb = $excel.ActiveWorkbook

# =========================================================================
# Adding a Tester profile
# =========================================================================

# ---- users sheet: add row 7 (Tester user) ----
$wsUsers = $wb.Worksheets.Item("users")
$wsUsers.Activate()
$wsUsers.Range("A7").Value = "test@test.com"
$wsUsers.Range("B7").Value = 12345678
$wsUsers.Range("C7").Value = "Tester"
$wsUsers.Range("D7").Value = "Testing"
$wsUsers.Range("G7").Value = $true
$wsUsers.Range("H7").Value = "admin"

# Make A7 a mailto hyperlink, matching the style used by the other
# hyperlinked emails in column A (copy format from A2).
$wsUsers.Hyperlinks.Add($wsUsers.Range("A7"), "mailto:test@test.com")
$wsUsers.Range("A2").Copy()
$wsUsers.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsUsers.Range("L7").Select()

# ---- watchlist sheet: add rows 18-21 for the new Tester (user_id 6) ----
$wsWatch = $wb.Worksheets.Item("watchlist")
$wsWatch.Activate()
$wsWatch.Range("A18").Value = 6
$wsWatch.Range("B18").Value = "TSLA"
$wsWatch.Range("A19").Value = 6
$wsWatch.Range("B19").Value = "AAPL"
$wsWatch.Range("A20").Value = 6
$wsWatch.Range("B20").Value = "MSFT"
$wsWatch.Range("A21").Value = 6
$wsWatch.Range("B21").Value = "DELL"
$wsWatch.Range("B21").Select()

# ---- portfolio sheet: add row 7 (Tester's Portfolio) ----
$wsPortfolio = $wb.Worksheets.Item("portfolio")
$wsPortfolio.Activate()
$wsPortfolio.Range("A7").Value = 6
$wsPortfolio.Range("B7").Value = "Tester's Portfolio"
$wsPortfolio.Range("C7").Value = 1000000
$wsPortfolio.Range("C9").Select()

# ---- portfolioprice sheet: add row 7 ----
$wsPrice = $wb.Worksheets.Item("portfolioprice")
$wsPrice.Activate()
$wsPrice.Range("A7").Value = 6
$wsPrice.Range("B7").Value = 6
$wsPrice.Range("C7").Value = 1000000
$wsPrice.Range("D7").Value = 0
$wsPrice.Range("E7").Value = 44120.275080960651
# Match the date/time number format used by the cell above it.
$wsPrice.Range("E6").Copy()
$wsPrice.Range("E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsPrice.Range("F7").Select()

# ---- transaction sheet: fill in rows 12-13 (previously blank placeholders) ----
$wsTxn = $wb.Worksheets.Item("transaction")
$wsTxn.Activate()
$wsTxn.Range("A12").Value = 6
$wsTxn.Range("B12").Value = 6
$wsTxn.Range("C12").Value = "DELL"
$wsTxn.Range("D12").Value = 60
$wsTxn.Range("E12").Value = 44119.483414351853
$wsTxn.Range("F12").Value = 2310
$wsTxn.Range("G12").Value = 0

$wsTxn.Range("A13").Value = 6
$wsTxn.Range("B13").Value = 6
$wsTxn.Range("C13").Value = "TSLA"
$wsTxn.Range("D13").Value = 60
$wsTxn.Range("E13").Value = 44119.483414351853
$wsTxn.Range("F13").Value = 3420
$wsTxn.Range("G13").Value = 0
$wsTxn.Range("D14").Select()

# ---- exchanges sheet: make it the active tab/sheet, as in the saved file ----
$wsExchanges = $wb.Worksheets.Item("exchanges")
$wsExchanges.Activate()
